$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '37.573.70'
$ws.Range('E2').Value = '  +0.27%  '
Set-TextValue 'D3' '2.083.20'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '235.41'
$ws.Range('E5').Value = '  +0.20%  '
Set-TextValue 'D6' '0.627'
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('E7').Value = '  -0.05%  '
Set-TextValue 'D8' '57.69'
$ws.Range('E8').Value = '  -1.04%  '
Set-TextValue 'D9' '0.389'
$ws.Range('E9').Value = '  +1.66%  '
Set-TextValue 'D10' '0.0780'
$ws.Range('E10').Value = '  +2.57%  '
$ws.Range('E11').Value = '  +1.32%  '
Set-TextValue 'D12' '2.390.56'
$ws.Range('E12').Value = '  +0.76%  '
Set-TextValue 'D13' '14.45'
$ws.Range('E13').Value = '  -0.51%  '
Set-TextValue 'D14' '20.86'
$ws.Range('E14').Value = '  -1.46%  '
Set-TextValue 'D15' '0.784'
$ws.Range('E15').Value = '  +0.80%  '
Set-TextValue 'D16' '5.24'
$ws.Range('E16').Value = '  +1.26%  '
Set-TextValue 'D17' '2.094.32'
$ws.Range('E17').Value = '  +1.30%  '
Set-TextValue 'D18' '37.525.54'
$ws.Range('E18').Value = '  -0.37%  '
Set-TextValue 'D19' '6.23'
$ws.Range('E19').Value = '  +0.72%  '
Set-TextValue 'D20' '69.77'
$ws.Range('E20').Value = '  -0.53%  '
Set-TextValue 'D21' '0.0₃0822'
$ws.Range('E21').Value = '  +0.79%  '
Set-TextValue 'D22' '226.94'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -0.03%  '
Set-TextValue 'D24' '2.43'
$ws.Range('E24').Value = '  +1.77%  '
Set-TextValue 'D25' '2.41'
$ws.Range('E25').Value = '  -1.20%  '
Set-TextValue 'D26' '168.82'
$ws.Range('E26').Value = '  +1.40%  '
Set-TextValue 'D27' '8.90'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('E28').Value = '  +4.30%  '
$ws.Range('E29').Value = '  -5.85%  '
Set-TextValue 'D30' '19.16'
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('E31').Value = '  -0.28%  '
Set-TextValue 'D32' '4.62'
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('E33').Value = '  -0.48%  '
Set-TextValue 'D34' '4.59'
$ws.Range('E34').Value = '  +0.30%  '
Set-TextValue 'D35' '2.52'
$ws.Range('E35').Value = '  -1.38%  '
Set-TextValue 'D36' '3.36'
$ws.Range('E36').Value = '  +0.18%  '
Set-TextValue 'D37' '1.79'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  +0.07%  '
Set-TextValue 'D39' '5.63'
$ws.Range('E39').Value = '  -3.78%  '
Set-TextValue 'D40' '2.94'
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('B41').Value = 'Cronos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D41' '0.0957'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D42' '1.489.76'
$ws.Range('E42').Value = '  +2.38%  '
Set-TextValue 'D43' '97.37'
$ws.Range('E43').Value = '  +1.66%  '
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('E45').Value = '  -1.38%  '
Set-TextValue 'D46' '4.20'
$ws.Range('E46').Value = '  -7.15%  '
$ws.Range('E47').Value = '  +0.89%  '
Set-TextValue 'D48' '15.57'
$ws.Range('E48').Value = '  -1.70%  '
Set-TextValue 'D49' '7.27'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('E50').Value = '  +1.23%  '
Set-TextValue 'D51' '2.278.05'
$ws.Range('E51').Value = '  +0.81%  '
